$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(11, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(11, 3).Value = "Coquimbo"
$ws.Cells.Item(11, 4).Value = 44463
$ws.Cells.Item(11, 5).Value = 4
$ws.Cells.Item(11, 6).Value = "Fruta"
$ws.Cells.Item(11, 7).Value = 100101
$ws.Cells.Item(11, 8).Value = "Berries"
$ws.Cells.Item(11, 9).Value = 100101001
$ws.Cells.Item(11, 10).Value = "Arándano (blue)"
$ws.Cells.Item(11, 11).Value = "Sin especificar"
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 100
$ws.Cells.Item(11, 14).Value = 13000
$ws.Cells.Item(11, 15).Value = 14000
$ws.Cells.Item(11, 16).Value = 13500
$ws.Cells.Item(11, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(11, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(11, 19).Value = 6750
$ws.Cells.Item(11, 20).Value = 2

$ws.Cells.Item(11, 4).NumberFormat = $ws.Cells.Item(10, 4).NumberFormat
